$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 "43.877.68"
Set-TextValue 2 5 "  +3.54%  "
Set-TextValue 3 4 "2.256.79"
Set-TextValue 3 5 "  +1.11%  "
Set-TextValue 4 5 "  +0.22%  "
Set-TextValue 5 4 "230.20"
Set-TextValue 5 5 "  -0.39%  "
Set-TextValue 6 4 "0.636"
Set-TextValue 6 5 "  +2.46%  "
Set-TextValue 7 4 "62.85"
Set-TextValue 7 5 "  +2.48%  "
Set-TextValue 8 5 "  +0.05%  "
Set-TextValue 9 4 "0.451"
Set-TextValue 9 5 "  +11.52%  "
Set-TextValue 10 5 "  +12.79%  "
Set-TextValue 11 4 "57.17"
Set-TextValue 11 5 "  -0.68%  "
Set-TextValue 12 5 "  +2.11%  "
Set-TextValue 13 4 "25.86"
Set-TextValue 13 5 "  +15.50%  "
Set-TextValue 14 4 "2.591.99"
Set-TextValue 14 5 "  +1.16%  "
Set-TextValue 15 4 "15.62"
Set-TextValue 15 5 "  +0.14%  "
Set-TextValue 16 4 "6.17"
Set-TextValue 16 5 "  +9.40%  "
Set-TextValue 17 4 "0.847"
Set-TextValue 17 5 "  +5.78%  "
Set-TextValue 18 4 "2.245.08"
Set-TextValue 18 5 "  -0.14%  "
Set-TextValue 19 4 "43.762.73"
Set-TextValue 19 5 "  +3.61%  "
Set-TextValue 20 5 "  +7.11%  "
Set-TextValue 21 4 "73.37"
Set-TextValue 21 5 "  +1.40%  "
Set-TextValue 22 4 "6.05"
Set-TextValue 22 5 "  -2.60%  "
Set-TextValue 23 4 "252.21"
Set-TextValue 23 5 "  +2.85%  "
Set-TextValue 24 5 "  +0.21%  "
Set-TextValue 25 4 "2.43"
Set-TextValue 25 5 "  +1.44%  "
Set-TextValue 26 5 "  -1.77%  "
Set-TextValue 27 4 "3.33"
Set-TextValue 27 5 "  +24.97%  "
Set-TextValue 28 4 "10.02"
Set-TextValue 28 5 "  +3.14%  "
Set-TextValue 29 4 "171.89"
Set-TextValue 29 5 "  +1.51%  "
Set-TextValue 30 5 "  -2.08%  "
Set-TextValue 31 4 "20.75"
Set-TextValue 31 5 "  +1.82%  "
Set-TextValue 32 5 "  -4.73%  "
Set-TextValue 33 5 "  +3.87%  "
Set-TextValue 34 5 "  +5.38%  "
Set-TextValue 35 4 "4.77"
Set-TextValue 35 5 "  +2.50%  "
Set-TextValue 36 4 "4.86"
Set-TextValue 36 5 "  -3.28%  "
Set-TextValue 37 5 "  +7.69%  "
Set-TextValue 38 4 "6.49"
Set-TextValue 38 5 "  +1.40%  "
Set-TextValue 39 4 "2.32"
Set-TextValue 39 5 "  -1.73%  "
Set-TextValue 40 4 "0.0257"
Set-TextValue 40 5 "  +2.56%  "
Set-TextValue 41 5 "  +0.19%  "
Set-TextValue 42 2 "TerraClassic"
Set-TextValue 42 3 "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue 42 4 "0.000219"
Set-TextValue 42 5 "  -4.86%  "
Set-TextValue 43 2 "InjectiveProtocol"
Set-TextValue 43 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 43 4 "17.52"
Set-TextValue 43 5 "  +8.38%  "
Set-TextValue 44 2 "Cronos"
Set-TextValue 44 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 44 4 "0.0974"
Set-TextValue 44 5 "  +1.05%  "
Set-TextValue 45 2 "FraxShare"
Set-TextValue 45 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 45 4 "8.26"
Set-TextValue 45 5 "  -4.36%  "
Set-TextValue 46 4 "98.00"
Set-TextValue 46 5 "  +0.78%  "
Set-TextValue 47 5 "  -1.36%  "
Set-TextValue 48 4 "4.33"
Set-TextValue 48 5 "  -0.37%  "
Set-TextValue 49 4 "1.445.33"
Set-TextValue 49 5 "  -1.15%  "
Set-TextValue 50 2 "NEARProtocol"
Set-TextValue 50 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 50 4 "2.30"
Set-TextValue 50 5 "  +2.81%  "
Set-TextValue 51 2 "Celestia"
Set-TextValue 51 3 "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue 51 4 "9.89"
Set-TextValue 51 5 "  +14.09%  "
